# Updates the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# with refreshed figures, matching the GitHub Actions data-refresh commit.
#
# Cells in column D whose new text parses as a plain number (e.g. "309.60",
# "0.000008720") are forced to the "@" (text) number format before the
# assignment so Excel keeps them as literal strings instead of silently
# coercing them to numeric values and dropping significant trailing zeros.
# Column E values are percentage strings (with surrounding spaces) and are
# never numeric, so they're assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.925.91"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "1.819.36"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.60"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4685"
$ws.Range("E7").Value = "  +3.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3696"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8727"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.49"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.815.47"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.369"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.67"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.510"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008720"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.78"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").Value = "26.950.72"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.346"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "2.041.51"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.902"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.44"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.188"
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.329"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08942"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7716"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.167"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.505"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.902"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.088"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01966"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05288"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.313"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5353"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.364"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1669"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.445"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4964"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.52"
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.671"
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.85"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  +0.20%  "
